$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column D (old column D "wait" shifts to E) -------------
# This mirrors the v0.4 test-case structure: a new "target" column is
# inserted between the existing "options" column (C) and the "wait" column.
$ws.Columns.Item(4).Insert()

# --- Row 1 (headers): duplicate the "assertAlert" header into the new D1 -
$ws.Cells.Item(1,4).Value = "assertAlert"

# --- Row 2 (the "open" test step): the old "options" payload that used to
# live in C3 (value/options-by-case) now lives in C2 (common options) -----
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(2,3).Value = "Hello World"

# --- Row 3 (the "assertAlert" test step): new D3 holds the "target"
# payload describing which case/selector the assertion applies to --------
$ws.Cells.Item(3,4).Value = '{"target":"Hello World"}'

# Give the new D3 cell its own look: a slightly larger green-ish font using
# the JP Gothic family (distinguishing "target" values from "options"
# values that keep the Sarasa Mono CL styling).
$f = $ws.Cells.Item(3,4).Font
$f.Size = 12
$f.Color = 1539334          # OLE BGR encoding of RGB(06,7D,17) -> FF067D17
$f.Name = "ＭＳ Ｐゴシック"

# --- Column widths ---------------------------------------------------------
# New column D needs to be wide enough to show the "target" JSON payload.
$ws.Columns.Item(4).ColumnWidth = 21.8

# --- Selection / cursor -----------------------------------------------------
$ws.Range("D9").Select()
